$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Shift-Cell($srcAddr, $dstAddr) {
  $src = $ws.Range($srcAddr)
  $dst = $ws.Range($dstAddr)
  [void]$src.Cut($dst)
}

# --- 1. Shift the K:P block one column right, to L:Q, preserving formatting ---
# Row 1: K1..P1 -> L1..Q1 (right to left so we never overwrite an unmoved source cell)
Shift-Cell "P1" "Q1"
Shift-Cell "O1" "P1"
Shift-Cell "N1" "O1"
Shift-Cell "M1" "N1"
Shift-Cell "L1" "M1"
Shift-Cell "K1" "L1"

# Row 2: K2..N2 -> L2..O2
Shift-Cell "N2" "O2"
Shift-Cell "M2" "N2"
Shift-Cell "L2" "M2"
Shift-Cell "K2" "L2"

# Row 3: K3..N3 -> L3..O3
Shift-Cell "N3" "O3"
Shift-Cell "M3" "N3"
Shift-Cell "L3" "M3"
Shift-Cell "K3" "L3"

# the old K column left an empty styled stub behind - drop it entirely
[void]$ws.Range("K1").Clear()

# --- 2. Re-capitalize the headers that changed text ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Quantity"
$ws.Range("D1").Value = "SL đã bán"
$ws.Range("E1").Value = "Số % giảm"
$ws.Range("F1").Value = "Giá"
$ws.Range("G1").Value = "Trạng thái"
$ws.Range("L1").Value = "Size "
$ws.Range("M1").Value = "Trọng lượng"
$ws.Range("N1").Value = "Trọng tải tối đa"
$ws.Range("O1").Value = "Chức năng"
$ws.Range("P1").Value = "Hình ảnh"
$ws.Range("Q1").Value = "Danh mục"

# --- 3. New "Loại sản phẩm" (product type) column J ---
$ws.Range("J1").Value = "Loại sản phẩm"
$ws.Range("J2").Value = 0
$ws.Range("J3").Value = 1

# --- 4. Column widths ---
$ws.Columns("A").ColumnWidth = 15
$ws.Columns("B").ColumnWidth = 21.1666666666667
$ws.Columns("C").ColumnWidth = 7.83333333333333
$ws.Columns("J").ColumnWidth = 13.1666666666667
$ws.Columns("Q").ColumnWidth = 13.1666666666667

# --- 5. Row height for the header row ---
$ws.Rows("1").RowHeight = 24

# --- 6. Selection matches the post-edit cursor position ---
[void]$ws.Range("P10").Select()
